$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Name"

$ws.Range("A5").Value = "Name"
$ws.Range("B5").Value = "Amount"
$ws.Range("C5").Value = "Price"

$ws.Range("A7").Value = "stood"
$ws.Range("B7").Value = "dood"
$ws.Range("C7").Value = "hood"

$ws.Range("A8").Value = "'626226"
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").Value = "'2626"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "226fdh"

$ws.Range("A9").Value = "'4444"
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").Value = "asdfdsf"
$ws.Range("C9").Value = "afdhfdh"

$ws.Range("A10").Value = "Name555"
$ws.Range("B10").Value = "Amount"
$ws.Range("C10").Value = "Price"
